# "Agregue para editar articulos" -- rework the "Hoja 1" product list so it
# can be edited: lowercase the existing headers, add "peso"/"stock" columns,
# refresh the "precio" values and fill in the new columns for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row --------------------------------------------------------
$ws.Range("B1").Value = "alimento"
$ws.Range("C1").Value = "tipo"
$ws.Range("D1").Value = "precio"
$ws.Range("E1").Value = "peso"
$ws.Range("F1").Value = "stock"

# --- Row 2: Nutribon / perro --------------------------------------------
$ws.Range("D2").Value = 2700
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 10

# --- Row 3: Dogui / perro ------------------------------------------------
$ws.Range("D3").Value = 9000
$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 10

# --- Row 4: Gatii / gato --------------------------------------------------
$ws.Range("D4").Value = 8900
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 20

# --- Stray formatting left on J11 (underline) + final selection ---------
$ws.Range("J11").Font.Underline = $true
$ws.Rows.Item(11).RowHeight = 15.75
[void]$ws.Range("J11").Select()

# --- Margins trimmed to zero --------------------------------------------
$ws.PageSetup.LeftMargin = 0
$ws.PageSetup.RightMargin = 0
$ws.PageSetup.TopMargin = 0
$ws.PageSetup.BottomMargin = 0
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0
